# Added Backup Process to StudentData
# Inserts a new backup/placeholder student record at row 51 of the
# "RawData" sheet, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 51..end down by one, inserting a blank row at 51.
$ws.Rows.Item(51).Insert()

# The freshly inserted row has no real data yet - strip any inherited
# formatting so the new row starts from the default style.
$ws.Range("A51:G51").ClearFormats()

# Populate the new backup record.
$ws.Range("A51").Value = 351353
$ws.Range("B51").Value = "the"
$ws.Range("C51").Value = "the"
$ws.Range("D51").Value = "efa"
$ws.Range("E51").Value = "Y"
$ws.Range("F51").Value = "N"
$ws.Range("G51").Value = "N"
